# Update "想去人数" (want-to-go count) values for a handful of exhibition
# rows that changed between scrapes, on both the "展览" sheet and the
# aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 3, 9, 10, 11 -> column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 532
$wsExhibit.Range("F9").Value = 331
$wsExhibit.Range("F10").Value = 3331
$wsExhibit.Range("F11").Value = 34

# Sheet "全部类型": rows 4, 10, 11, 12 -> column F (same events, shifted by one row)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 532
$wsAll.Range("F10").Value = 331
$wsAll.Range("F11").Value = 3331
$wsAll.Range("F12").Value = 34
